$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Responsible" column (B) with names, right-aligned
$ws.Range("B4").Value = "Aziz"
$ws.Range("B5").Value = "Aziz"
$ws.Range("B3").Value = "Mark"

$ws.Range("B3:C5").HorizontalAlignment = -4152  # xlRight

# Update the active selection to reflect where the user ended up
$ws.Range("F11").Select()
